$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 366967
$ws.Range("I9").Value = 366967
$ws.Range("K9").Value = 366967
$ws.Range("M9").Value = -366798
$ws.Range("H17").Value = 2810.125
$ws.Range("J17").Value = 2810.125
$ws.Range("L17").Value = 8430.375
$ws.Range("N17").Value = -8766.375
$ws.Range("H19").Value = 923.6
$ws.Range("I19").Value = 810.8333
$ws.Range("J19").Value = 1092.75
$ws.Range("K19").Value = 810.8333
$ws.Range("L19").Value = 1092.75
$ws.Range("M19").Value = -635.8333
$ws.Range("N19").Value = -1442.75
$ws.Range("H28").Value = 11907213
$ws.Range("I28").Value = 15153303
$ws.Range("J28").Value = 4880.1665
$ws.Range("K28").Value = 15153303
$ws.Range("L28").Value = 4880.1665
$ws.Range("M28").Value = -15152818
$ws.Range("N28").Value = -5850.1665
$ws.Range("H62").Value = 1824
$ws.Range("I62").Value = 1824
$ws.Range("K62").Value = 1824
$ws.Range("M62").Value = -1200
$ws.Range("H65").Value = 1824
$ws.Range("I65").Value = 1824
$ws.Range("K65").Value = 9120
$ws.Range("M65").Value = -6000
$ws.Range("H98").Value = 7925.579
$ws.Range("I98").Value = 1545.1538
$ws.Range("J98").Value = 21749.834
$ws.Range("K98").Value = 1545.1538
$ws.Range("L98").Value = 21749.834
$ws.Range("M98").Value = -47.15380000000005
$ws.Range("N98").Value = -24745.834
$ws.Range("H112").Value = 6457.9707
$ws.Range("I112").Value = 4100
$ws.Range("J112").Value = 6772.3667
$ws.Range("K112").Value = 12300
$ws.Range("L112").Value = 20317.1001
$ws.Range("M112").Value = -11192
$ws.Range("N112").Value = -22533.1001
$ws.Range("H116").Value = 4999.3335
$ws.Range("I116").Value = 4999.3335
$ws.Range("K116").Value = 4999.3335
$ws.Range("M116").Value = -1557.3335
$ws.Range("H122").Value = 7925.579
$ws.Range("I122").Value = 1545.1538
$ws.Range("J122").Value = 21749.834
$ws.Range("K122").Value = 4635.4614
$ws.Range("L122").Value = 65249.50199999999
$ws.Range("M122").Value = -2185.4614
$ws.Range("N122").Value = -70149.50199999999
$ws.Range("H131").Value = 3539
$ws.Range("I131").Value = 1916
$ws.Range("K131").Value = 5748
$ws.Range("M131").Value = -708
$ws.Range("H135").Value = 3364.75
$ws.Range("I135").Value = 3182.8
$ws.Range("K135").Value = 28645.2
$ws.Range("M135").Value = -26110.2
$ws.Range("H138").Value = 3366.1553
$ws.Range("J138").Value = 2747.2666
$ws.Range("L138").Value = 8241.799800000001
$ws.Range("N138").Value = -18521.7998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2113.5715
$ws.Range("I45").Value = 1674.375
$ws.Range("J45").Value = 3519
$ws.Range("K45").Value = 1674.375
$ws.Range("L45").Value = 3519
$ws.Range("M45").Value = -1297.375
$ws.Range("N45").Value = -4273
$ws.Range("H122").Value = 1515.138
$ws.Range("I122").Value = 768.5
$ws.Range("K122").Value = 2305.5
$ws.Range("M122").Value = 144.5
$ws.Range("H132").Value = 2119.12
$ws.Range("I132").Value = 1989.5238
$ws.Range("J132").Value = 2799.5
$ws.Range("K132").Value = 5968.5714
$ws.Range("L132").Value = 8398.5
$ws.Range("M132").Value = -3438.5714
$ws.Range("N132").Value = -13458.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 98000
$ws.Range("J132").Value = 98000
$ws.Range("L132").Value = 98000
$ws.Range("N132").Value = -108120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 11996.667
$ws.Range("I51").Value = 11996.667
$ws.Range("K51").Value = 11996.667
$ws.Range("M51").Value = -11260.667
$ws.Range("H61").Value = 11996.667
$ws.Range("I61").Value = 11996.667
$ws.Range("K61").Value = 11996.667
$ws.Range("M61").Value = -11648.667
$ws.Range("H99").Value = 1771.5625
$ws.Range("I99").Value = 1771.5625
$ws.Range("K99").Value = 1771.5625
$ws.Range("M99").Value = -273.5625
$ws.Range("H105").Value = 3607.4167
$ws.Range("I105").Value = 3258.9
$ws.Range("K105").Value = 3258.9
$ws.Range("M105").Value = -1511.9
$ws.Range("H106").Value = 90000
$ws.Range("J106").Value = 90000
$ws.Range("L106").Value = 90000
$ws.Range("N106").Value = -92524
$ws.Range("H126").Value = 1771.5625
$ws.Range("I126").Value = 1771.5625
$ws.Range("K126").Value = 5314.6875
$ws.Range("M126").Value = -2844.6875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 18498.5
$ws.Range("I74").Value = 15495.5
$ws.Range("J74").Value = 20000
$ws.Range("K74").Value = 46486.5
$ws.Range("L74").Value = 60000
$ws.Range("M74").Value = -45425.5
$ws.Range("N74").Value = -62122
$ws.Range("H77").Value = 18498.5
$ws.Range("I77").Value = 15495.5
$ws.Range("J77").Value = 20000
$ws.Range("K77").Value = 139459.5
$ws.Range("L77").Value = 180000
$ws.Range("M77").Value = -134155.5
$ws.Range("N77").Value = -190608
$ws.Range("H137").Value = 67883.234
$ws.Range("I137").Value = 150506.42
$ws.Range("J137").Value = 10047
$ws.Range("K137").Value = 451519.26
$ws.Range("L137").Value = 30141
$ws.Range("M137").Value = -446419.26
$ws.Range("N137").Value = -40341
$ws.Range("H138").Value = 16517.523
$ws.Range("J138").Value = 21072.467
$ws.Range("L138").Value = 63217.401
$ws.Range("N138").Value = -73497.401

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 27214.4
$ws.Range("J98").Value = 27214.4
$ws.Range("L98").Value = 27214.4
$ws.Range("N98").Value = -33204.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 579.1111
$ws.Range("I16").Value = 563.75
$ws.Range("J16").Value = 702
$ws.Range("K16").Value = 563.75
$ws.Range("L16").Value = 702
$ws.Range("M16").Value = -393.75
$ws.Range("N16").Value = -1042
$ws.Range("H42").Value = 20000
$ws.Range("I42").Value = 20000
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 20000
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -19437
$ws.Range("N42").ClearContents()
$ws.Range("H46").Value = 2063.7334
$ws.Range("I46").Value = 1574
$ws.Range("J46").Value = 2241.818
$ws.Range("K46").Value = 1574
$ws.Range("L46").Value = 2241.818
$ws.Range("M46").Value = -1386
$ws.Range("N46").Value = -2617.818
$ws.Range("H49").Value = 20000
$ws.Range("I49").Value = 20000
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 20000
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -19853
$ws.Range("N49").ClearContents()
$ws.Range("H68").Value = 2498.5
$ws.Range("I68").Value = 2481.3333
$ws.Range("J68").Value = 2550
$ws.Range("K68").Value = 2481.3333
$ws.Range("L68").Value = 2550
$ws.Range("M68").Value = -1732.3333
$ws.Range("N68").Value = -4048
$ws.Range("H71").Value = 2498.5
$ws.Range("I71").Value = 2481.3333
$ws.Range("J71").Value = 2550
$ws.Range("K71").Value = 12406.6665
$ws.Range("L71").Value = 12750
$ws.Range("M71").Value = -8662.666499999999
$ws.Range("N71").Value = -20238
$ws.Range("H128").Value = 84999.5
$ws.Range("J128").Value = 84999.5
$ws.Range("L128").Value = 84999.5
$ws.Range("N128").Value = -94959.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4474.316
$ws.Range("I122").Value = 1960.2059
$ws.Range("K122").Value = 5880.6177
$ws.Range("M122").Value = -3430.6177
$ws.Range("H136").Value = 3755.4517
$ws.Range("I136").Value = 3549.577
$ws.Range("J136").Value = 4826
$ws.Range("K136").Value = 10648.731
$ws.Range("L136").Value = 14478
$ws.Range("M136").Value = -8098.731
$ws.Range("N136").Value = -19578
